$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45177 -> 45178) for every data row (rows 2 through 468).
$lastRow = 468
$ws.Range("C2:C$lastRow").Value = 45178
